$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1760797342192691
$ws.Range("C2").Value = 0.5946843853820598
$ws.Range("J2").Value = 0.0132890365448505
$ws.Range("P2").Value = 0.1262458471760797
$ws.Range("S2").Value = 0.08970099667774087
$ws.Range("B3").Value = 0.01047120418848168
$ws.Range("C3").Value = 0.04712041884816754
$ws.Range("J3").Value = 0.02617801047120419
$ws.Range("P3").Value = 0.7591623036649214
$ws.Range("S3").Value = 0.1570680628272251
$ws.Range("J4").Value = 0.02380952380952381
$ws.Range("P4").Value = 0.6904761904761905
$ws.Range("S4").Value = 0.2857142857142857
$ws.Range("B6").Value = 0.06278026905829596
$ws.Range("D6").Value = 0.004484304932735426
$ws.Range("F6").Value = 0.04035874439461883
$ws.Range("J6").Value = 0.2197309417040359
$ws.Range("O6").Value = 0.02242152466367713
$ws.Range("Q6").Value = 0.1300448430493273
$ws.Range("R6").Value = 0.07174887892376682
$ws.Range("S6").Value = 0.4484304932735426
$ws.Range("B7").Value = 0.07909604519774012
$ws.Range("F7").Value = 0.06779661016949153
$ws.Range("J7").Value = 0.1186440677966102
$ws.Range("O7").Value = 0.01129943502824859
$ws.Range("Q7").Value = 0.1977401129943503
$ws.Range("R7").Value = 0.1016949152542373
$ws.Range("S7").Value = 0.423728813559322
$ws.Range("B8").Value = 0.1089965397923875
$ws.Range("D8").Value = 0.01557093425605536
$ws.Range("E8").Value = 0.001730103806228374
$ws.Range("F8").Value = 0.05882352941176471
$ws.Range("J8").Value = 0.1072664359861592
$ws.Range("O8").Value = 0.01903114186851211
$ws.Range("Q8").Value = 0.1782006920415225
$ws.Range("R8").Value = 0.08996539792387544
$ws.Range("S8").Value = 0.4204152249134948
$ws.Range("B9").Value = 0.106951871657754
$ws.Range("D9").Value = 0.0053475935828877
$ws.Range("F9").Value = 0.0748663101604278
$ws.Range("J9").Value = 0.1016042780748663
$ws.Range("O9").Value = 0.0267379679144385
$ws.Range("Q9").Value = 0.1390374331550802
$ws.Range("R9").Value = 0.06951871657754011
$ws.Range("S9").Value = 0.4759358288770054
$ws.Range("B10").Value = 0.09495982468955443
$ws.Range("D10").Value = 0.02337472607742878
$ws.Range("F10").Value = 0.06647187728268809
$ws.Range("J10").Value = 0.1176040905770635
$ws.Range("O10").Value = 0.01314828341855369
$ws.Range("Q10").Value = 0.2176771365960555
$ws.Range("R10").Value = 0.07815924032140248
$ws.Range("S10").Value = 0.3886048210372535
$ws.Range("G11").Value = 0.1449275362318841
$ws.Range("J11").Value = 0.06159420289855073
$ws.Range("K11").Value = 0.1847826086956522
$ws.Range("L11").Value = 0.5797101449275363
$ws.Range("S11").Value = 0.02898550724637681
$ws.Range("G12").Value = 0.7735849056603774
$ws.Range("J12").Value = 0.169811320754717
$ws.Range("K12").Value = 0.01257861635220126
$ws.Range("L12").Value = 0.006289308176100629
$ws.Range("S12").Value = 0.03773584905660377
$ws.Range("G13").Value = 0.6
$ws.Range("J13").Value = 0.3666666666666666
$ws.Range("S13").Value = 0.03333333333333333
$ws.Range("F15").Value = 0.02489626556016597
$ws.Range("H15").Value = 0.1576763485477178
$ws.Range("I15").Value = 0.07468879668049792
$ws.Range("J15").Value = 0.3858921161825726
$ws.Range("K15").Value = 0.04564315352697095
$ws.Range("M15").Value = 0.01659751037344398
$ws.Range("O15").Value = 0.07468879668049792
$ws.Range("S15").Value = 0.2199170124481328
$ws.Range("F16").Value = 0.01477832512315271
$ws.Range("H16").Value = 0.1773399014778325
$ws.Range("I16").Value = 0.07881773399014778
$ws.Range("J16").Value = 0.4187192118226601
$ws.Range("K16").Value = 0.07389162561576355
$ws.Range("M16").Value = 0.004926108374384237
$ws.Range("O16").Value = 0.08374384236453201
$ws.Range("S16").Value = 0.1477832512315271
$ws.Range("F17").Value = 0.02277432712215321
$ws.Range("H17").Value = 0.2173913043478261
$ws.Range("I17").Value = 0.08074534161490683
$ws.Range("J17").Value = 0.4492753623188406
$ws.Range("K17").Value = 0.06625258799171843
$ws.Range("M17").Value = 0.006211180124223602
$ws.Range("N17").Value = 0.004140786749482402
$ws.Range("O17").Value = 0.06004140786749482
$ws.Range("S17").Value = 0.09316770186335403
$ws.Range("H18").Value = 0.2178217821782178
$ws.Range("I18").Value = 0.07425742574257425
$ws.Range("J18").Value = 0.4653465346534654
$ws.Range("K18").Value = 0.06930693069306931
$ws.Range("M18").Value = 0.004950495049504951
$ws.Range("O18").Value = 0.07425742574257425
$ws.Range("S18").Value = 0.09405940594059406
$ws.Range("F19").Value = 0.007735583684950774
$ws.Range("H19").Value = 0.2531645569620253
$ws.Range("I19").Value = 0.06962025316455696
$ws.Range("J19").Value = 0.3663853727144866
$ws.Range("K19").Value = 0.1054852320675106
$ws.Range("M19").Value = 0.01617440225035162
$ws.Range("N19").Value = 0.0007032348804500703
$ws.Range("O19").Value = 0.06680731364275667
$ws.Range("S19").Value = 0.1139240506329114
